# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the slide master (rId12 in
#                             slideMaster1.xml.rels) and currently holds the
#                             "Integral" color scheme.
#   ppt/theme/theme2.xml  -> bound to the notes master and currently holds
#                             the default "Office Theme" color scheme.
# The target edit swaps the two color schemes: the slide master's theme
# becomes the default Office colors, and the notes master's theme becomes
# the Integral colors. The font scheme and format scheme (fills/lines/
# effects) are already byte-identical between the two theme parts, so the
# only real content delta is the 12-color scheme (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink).
#
# Apply the new ("Office Theme") color values to the presentation's theme
# via the slide master's Theme.ThemeColorScheme, which is the supported,
# persisted path for recoloring a theme in this object model.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Office Theme default color scheme (RGB packed as 0xBBGGRR, matching the
# PowerPoint COM RGB() convention) in clrScheme order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6,
# hlink, folHlink
$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
